# Insert a new data row at row 181 (pushing existing rows 181-289 down to
# 182-290) and populate it with the new reading, matching the target diff:
#   - dimension grows from A1:R289 to A1:R290
#   - new row 181 holds a fresh "Albahaca" observation for Región
#     Metropolitana dated 44582 ($/docena de matas, 542 $/Kg, 6 Kg o Unid.)
#   - every subsequent row simply shifts down by one (old row181 -> row182,
#     ..., old row289 -> row290)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 181:289 down to 182:290, leaving row 181 blank (inherits the
# date-style formatting already applied to column D below it, same as
# Excel's native "Insert Copied/Entire Row" behaviour).
$ws.Rows("181:181").Insert()

# Populate the newly-inserted row 181 with the new observation.
$ws.Cells.Item(181, 1).Value = 9
$ws.Cells.Item(181, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(181, 3).Value = "Metropolitana"
$ws.Cells.Item(181, 4).Value = 44582
$ws.Cells.Item(181, 5).Value = 13
$ws.Cells.Item(181, 6).Value = 100112052
$ws.Cells.Item(181, 7).Value = "Albahaca"
$ws.Cells.Item(181, 8).Value = "Sin especificar"
$ws.Cells.Item(181, 9).Value = "Primera"
$ws.Cells.Item(181, 10).Value = 250
$ws.Cells.Item(181, 11).Value = 3000
$ws.Cells.Item(181, 12).Value = 3500
$ws.Cells.Item(181, 13).Value = 3250
$ws.Cells.Item(181, 14).Value = "`$/docena de matas"
$ws.Cells.Item(181, 15).Value = "Región Metropolitana"
$ws.Cells.Item(181, 16).Value = 542
$ws.Cells.Item(181, 17).Value = 6
$ws.Cells.Item(181, 18).Value = "Hortaliza"
